# Weekly fruit/vegetable price update: insert a new "Haba" (Fava bean) price
# record for Femacal de La Calera / Coquimbo just before the current row 241,
# shifting the existing rows 241-276 down to 242-277.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 241..276 down to 242..277 and leave a fresh blank row at 241.
$ws.Rows.Item(241).Insert()

# Populate the newly inserted row 241 with the new weekly observation.
$ws.Cells.Item(241, 1).Value = 3
$ws.Cells.Item(241, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(241, 3).Value = "Coquimbo"
$ws.Cells.Item(241, 4).Value = 45127
$ws.Cells.Item(241, 5).Value = 5
$ws.Cells.Item(241, 6).Value = 100112026
$ws.Cells.Item(241, 7).Value = "Haba"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 75
$ws.Cells.Item(241, 11).Value = 15000
$ws.Cells.Item(241, 12).Value = 16000
$ws.Cells.Item(241, 13).Value = 15467
$ws.Cells.Item(241, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(241, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(241, 16).Value = 619
$ws.Cells.Item(241, 17).Value = 25
$ws.Cells.Item(241, 18).Value = "Hortaliza"
